# Daily refresh of the water-delivery tracker sheet.
#
# Columns: D = total cycle length in days, E = remaining days left in the
# current cycle, F = the date (yyyyMMdd) the current cycle started.
#
# Every day the sheet is refreshed against "today": the remaining-days
# counter is recomputed as (total - days elapsed since the cycle start).
# Once a cycle is exhausted (elapsed >= total, i.e. remaining would drop to
# 0 or below) it rolls straight into a fresh cycle: remaining resets to the
# full total and the start date becomes today.
#
# Rows whose start date isn't a parseable yyyyMMdd value (bad/legacy data)
# are left untouched rather than guessed at.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$todayStr = "20260126"
$todaySerial = [datetime]::ParseExact($todayStr, "yyyyMMdd", $null).ToOADate()

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {

    $totalCell = $ws.Cells.Item($row, 4)   # D
    $startCell = $ws.Cells.Item($row, 6)   # F

    $total = $totalCell.Value2
    if ($null -eq $total) { continue }

    $startRaw = [string]$startCell.Value2

    try {
        $startSerial = [datetime]::ParseExact($startRaw, "yyyyMMdd", $null).ToOADate()
    } catch {
        # Not a real date (corrupted/legacy row) -> skip, leave as-is.
        continue
    }

    $elapsedDays = [int]($todaySerial - $startSerial)
    $remaining = [int]$total - $elapsedDays

    if ($remaining -le 0) {
        # Cycle finished -> start a brand-new one today.
        $ws.Cells.Item($row, 5).Value = [int]$total
        $ws.Cells.Item($row, 6).Value = [int]$todayStr
    } else {
        $ws.Cells.Item($row, 5).Value = $remaining
    }
}
